# Generate Report for Handback
# Updates the handback-status workbook to reflect that the
# d6875cf9-4019-4dc6-b293-0d94272724f5 file has now been processed:
#  - Overview sheet: refresh "Latest HO Xliff Generate Date" for that row
#  - zh-cn / de-de sheets: mark "Content Duplicate" True and stamp the
#    real Correspond Handoff / Handback datetimes (previously just copied
#    from the other row as placeholders).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-19 18:56:12"

# ---- zh-cn sheet ----
$zhcn = $wb.Worksheets.Item("zh-cn")

# Content Duplicate (row 3) changes from False to True. Use
# SetCellDataTypeFromCell against a cell that already holds the text
# "True" so Excel stores it as a shared string rather than auto-coercing
# it into a boolean cell.
$zhcn.Range("F3").SetCellDataTypeFromCell($zhcn.Range("M2"), "True")

$zhcn.Range("H3").Value = "2016-08-19 18:56:01"
$zhcn.Range("K3").Value = "2016-08-19 18:56:29"

# ---- de-de sheet ----
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("F3").SetCellDataTypeFromCell($dede.Range("M2"), "True")

$dede.Range("H3").Value = "2016-08-19 18:56:12"
$dede.Range("K3").Value = "2016-08-19 18:56:36"
